$d = $word.ActiveDocument

$d.Content.Find.Execute("Alternate", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Opt", 2)

$d.Content.Find.Execute("4", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3", 2)
